$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "DRV8825 Stepper Motor Driver" row entirely; everything below
# shifts up by one (formulas, hyperlinks, formatting move with it).
$ws.Rows(5).Delete()

# The amount ordered for the 8-Channel Relay changed from 1 to 0 (it is now
# at row 11 after the deletion above).
$ws.Range("C11").Value = 0

# Insert a new row just above the table's totals ("Sum") row so the new
# part doesn't clobber the totals row, and so the table/totals machinery
# keeps working.
$ws.Rows(18).Insert()

# New part: "Other stuff"
$ws.Range("B18").Value = "Other stuff"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 25
$ws.Range("E18").Formula = "=C18*D18"

# Make sure the table definition covers the newly inserted row plus the
# (now shifted down) totals row again.
$tbl = $ws.ListObjects("Tabelle1")
$tbl.Resize($ws.Range("B4:E19"))

# New summary line below the table: price per person (split 4 ways).
$ws.Range("B21").Value = "Price per Person"
$ws.Range("B21").Font.Bold = $true

$ws.Range("E21").Formula = "=Tabelle1[[#Totals],[Price]] / 4"
$ws.Range("E21").Font.Bold = $true
$ws.Range("E21").NumberFormat = "#,##0.00\ ""€"""

# Update the remembered selection to match the authored workbook.
$ws.Range("D22").Select()

# Force a full recalculation so the cached formula results (Sum, Price per
# Person, ...) are refreshed even though the workbook is in manual
# calculation mode.
$excel.CalculateFull()
